# Thêm alert check trường thông tin nhập toa thuốc và sửa các lastIndex nếu danh sách rỗng
#
# This script adds a new prescription ("Toa thuốc cho Riêu") together with its
# medicine line (THUOC_TRONG_TOA) and also fixes up the names of the two
# existing prescriptions, matching the data entered through the application's
# "add prescription" form.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PRESCRIPTION sheet
# ---------------------------------------------------------------------------
$presSheet = $wb.Worksheets.Item("PRESCRIPTION")

# Fix the name of the first two existing prescriptions
$presSheet.Range("B2").Value = "Toa thuốc số 1"
$presSheet.Range("B3").Value = "Toa thuốc chữa trĩ cho Hồng"

# Update the start/end dates of the first two prescriptions
$presSheet.Range("C2").Value = 44757
$presSheet.Range("D2").Value = 44755
$presSheet.Range("C3").Value = 44751
$presSheet.Range("D3").Value = 44745
$presSheet.Range("C2:D3").NumberFormat = "dd/mm/yyyy"

# Add the new, third prescription
$presSheet.Range("A4").Value = 3
$presSheet.Range("B4").Value = "Toa thuốc cho Riêu"
$presSheet.Range("C4").Value = 44761
$presSheet.Range("D4").Value = 44757
$presSheet.Range("C4:D4").NumberFormat = "dd/mm/yyyy"

# ---------------------------------------------------------------------------
# THUOC_TRONG_TOA sheet (medicine lines belonging to each prescription)
# ---------------------------------------------------------------------------
$lineSheet = $wb.Worksheets.Item("THUOC_TRONG_TOA")

# Existing rows now point at the corrected medicine/unit/dosage for
# prescriptions 1 and 2
$lineSheet.Range("B2").Value = "Acemol"
$lineSheet.Range("C2").Value = "Hop"
$lineSheet.Range("D2").Value = "1 viên 1 ngày sau ăn"

$lineSheet.Range("B3").Value = "Avarino"
$lineSheet.Range("C3").Value = "Vi"
$lineSheet.Range("D3").Value = "1 viên 1 ngày"

# Add the medicine line for the new, third prescription
$lineSheet.Range("A4").Value = 1
$lineSheet.Range("B4").Value = "Ausagel "
$lineSheet.Range("C4").Value = "Vi"
$lineSheet.Range("D4").Value = "150ml 1 ngày"
$lineSheet.Range("E4").Value = 3
